$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: "3 months" header, matching the bold/bordered header style
# used by the rest of row 1 (copy format from the neighboring header cell).
$ws.Range("K1").Value = "3 months"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data row for the latest run
$ws.Range("A9").Value = 0.00029999999999999997
$ws.Range("B9").Value = 0.00088951110248452342
$ws.Range("C9").Value = 0.2129999999999998
$ws.Range("D9").Value = 64
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = "set to 5"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0

$ws.Range("L14").Select()
